# Applies the "IMAGING_DATASET" -> "IMAGING_DATASET_ID" header rename on both
# sheets, plus fills in the previously-empty INDIVIDUAL_REF_ID / INDIVIDUAL_REF_DB
# (columns BL/BM) example & guide rows on the "Examples & Info" sheet.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("Examples & Info")

# --- Rename header BA1 on both sheets ---
$wsData.Range("BA1").Value = "IMAGING_DATASET_ID"
$wsInfo.Range("BA1").Value = "IMAGING_DATASET_ID"

# --- Fill in the new INDIVIDUAL_REF_ID (BL) / INDIVIDUAL_REF_DB (BM) columns ---
# Row 2: column description text
$wsInfo.Range("BL2").Value = "ID/Label of the individual as referenced to in external database"
$wsInfo.Range("BM2").Value = "Label of the external mouse database e.g. MoVi / CRUK-CI"

# Row 3: example value
$wsInfo.Range("BL3").Value = "[12345]"
$wsInfo.Range("BM3").Value = "['MoVi, CRUK-CI']"

# Row 4: lab/guide value
$wsInfo.Range("BL4").Value = "Odomlab"
$wsInfo.Range("BM4").Value = "Odomlab"

# Row 5: category value
$wsInfo.Range("BL5").Value = "sample"
$wsInfo.Range("BM5").Value = "sample"
